# [Build 401] [Update] login et compte commence
# Updates the "SplashScreen" profile-creation rows and inserts a new
# "Ajouter l'ecran de ranking" task row in the Planning worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 20: SplashScreen / (new) Team Blui ------------------------------
$ws.Range("C20").Value = "Team Blui"
$ws.Range("E20").Value = 0.2
$ws.Range("F20").Value = 'Le "cœur" a été desinné = A améliorer. Faire le cube team blui qui passe et qui s''écrase avec d''autres cubes'

# --- Row 21: SplashScreen / Init de profile / Retrouver un profil pas save
# Copy the shaded ("green") formatting used elsewhere in the table (e.g.
# row 16) onto B21:D21, and the shaded percent formatting onto E21.
$ws.Range("B16:D16").Copy()
$ws.Range("B21:D21").PasteSpecial(-4122)
$ws.Range("E16").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B21").Value = "SplashScreen"
$ws.Range("C21").Value = "Init de profile"
$ws.Range("D21").Value = "Retrouver un profil pas save"
$ws.Range("E21").Value = 1

# --- Row 22: SplashScreen / Nouveau profil --------------------------------
$ws.Range("C22").Value = "Nouveau profil"

# --- Row 23: SplashScreen / Profil existant -------------------------------
$ws.Range("C23").Value = "Profil existant"

# --- Insert new row 52: Wheelsong / Ajouter l'ecran de ranking -----------
$ws.Rows("52:52").Insert()
$ws.Range("B52").Value = "Wheelsong"
$ws.Range("D52").Value = "Ajouter l'écran de ranking"
$ws.Range("E52").Value = 0

# --- Refresh the on-screen selection to match the author's final cursor --
$ws.Range("F21").Select() | Out-Null
